$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 44-45 (pushes the former rows 44-54 down to 46-56).
$ws.Range("A44:A45").EntireRow.Insert()

# Row 44 - new weekly data entry
$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C44").Value = "Ñuble"
$ws.Range("D44").Value = 45205
$ws.Range("E44").Value = 16
$ws.Range("F44").Value = 300000000
$ws.Range("G44").Value = "Espárragos"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 400
$ws.Range("K44").Value = 1300
$ws.Range("L44").Value = 1300
$ws.Range("M44").Value = 1300
$ws.Range("N44").Value = "`$/kilo"
$ws.Range("O44").Value = "Provincia de Diguillín"
$ws.Range("P44").Value = 1300
$ws.Range("Q44").Value = 1
$ws.Range("R44").Value = "Hortaliza"

# Row 45 - new weekly data entry
$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C45").Value = "Ñuble"
$ws.Range("D45").Value = 45205
$ws.Range("E45").Value = 16
$ws.Range("F45").Value = 300000000
$ws.Range("G45").Value = "Espárragos"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 300
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 1500
$ws.Range("M45").Value = 1500
$ws.Range("N45").Value = "`$/kilo"
$ws.Range("O45").Value = "Región del Maule"
$ws.Range("P45").Value = 1500
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = "Hortaliza"
